$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 408
$ws1.Range("F5").Value = 1108
$ws1.Range("F8").Value = 835
$ws1.Range("F9").Value = 1588
$ws1.Range("G9").Value = 61.2
$ws1.Range("F10").Value = 6018
$ws1.Range("F12").Value = 1718
$ws1.Range("F13").Value = 434
$ws1.Range("F14").Value = 5802
$ws1.Range("F15").Value = 110
$ws1.Range("F17").Value = 148
$ws1.Range("F19").Value = 1624
$ws1.Range("F22").Value = 136
$ws1.Range("F23").Value = 1313
$ws1.Range("F24").Value = 712
$ws1.Range("F30").Value = 3848

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 151
$ws2.Range("F8").Value = 368

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2212
$ws3.Range("F4").Value = 591
$ws3.Range("F5").Value = 154

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2212
$ws4.Range("F4").Value = 591
$ws4.Range("F6").Value = 408
$ws4.Range("F7").Value = 1108
$ws4.Range("F12").Value = 835
$ws4.Range("F13").Value = 154
$ws4.Range("F14").Value = 1588
$ws4.Range("G14").Value = 61.2
$ws4.Range("F15").Value = 6018
$ws4.Range("F17").Value = 1718
$ws4.Range("F20").Value = 434
$ws4.Range("F23").Value = 5802
$ws4.Range("F24").Value = 110
$ws4.Range("F26").Value = 148
$ws4.Range("F28").Value = 1624
$ws4.Range("F31").Value = 136
$ws4.Range("F32").Value = 1313
$ws4.Range("F33").Value = 712
$ws4.Range("F45").Value = 3848
